# Insert an "Application" value into column C for rows 2-20, shifting the
# existing contents of columns C:F one column to the right (C->D, D->E, E->F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cCell = $ws.Cells.Item($row, 3)  # C
    $dCell = $ws.Cells.Item($row, 4)  # D
    $eCell = $ws.Cells.Item($row, 5)  # E
    $fCell = $ws.Cells.Item($row, 6)  # F

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $eVal = $eCell.Value()

    $fCell.Value = $eVal
    $eCell.Value = $dVal
    $dCell.Value = $cVal
    $cCell.Value = "Application"
}
